$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -152.0046
$ws.Range("B2").Value = -151.9461

$ws.Range("A3").Value = 59.2679
$ws.Range("B3").Value = 59.3039

$ws.Range("A4").Value = -150.7751
$ws.Range("B4").Value = -150.834

$ws.Range("A5").Value = 59.8943
$ws.Range("B5").Value = 59.8581
